# The KiCad-generated BOM sheet was regenerated: the resistor group in row 6
# lost designator R2 (now R3,R8,R9,R5,R4) and its quantity dropped from 6 to 5.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "R3,R8,R9,R5,R4"

# Restore the last-used selection/cursor cell recorded in the saved view state.
$ws.Range("B7").Select()
